$d = $word.ActiveDocument

# ------------------------------------------------------------------
# Goal (per the diff):
#   "...lignment Metamodels and Big Data)."
# becomes three runs:
#   1) "...lignment Metamodels and Big Data. "   (same Arial formatting as before)
#   2) "'Distance' calculation models between two aligned resources in a
#       given 'axis' or parent class"            (plain / rtl=0 formatting)
#   3) ")."                                      (same Arial formatting as run 1)
# ------------------------------------------------------------------

# Step 1: rewrite the existing run's text, inserting a space after "Data."
# and keeping the trailing "). " for now - this all stays inside the single
# original (Arial-formatted) run, so its formatting is preserved exactly.
$target = $d.Content
$target.Find.Execute(
    "lignment Metamodels and Big Data).", $true, $false, $false, $false,
    $false, $true, 1, $false, "lignment Metamodels and Big Data. ).", 2
) | Out-Null

$runStart = $target.Start

# Step 2: build the new sentence in a throw-away paragraph at the end of the
# document. A brand new paragraph's run only ever picks up the minimal
# "rtl=0" run properties (no explicit font overrides), which is exactly the
# formatting the diff shows for the inserted sentence.
$scratchParagraph = $d.Paragraphs.Add()
$scratchParagraph.Range.Text = "‘Distance’ calculation models between two aligned resources in a given ‘axis’ or parent class"

$scratchRange = $scratchParagraph.Range
$scratchRange.MoveEnd(1, -1) | Out-Null
Write-Output "scratch text to copy: '$($scratchRange.Text)'"
$scratchRange.Copy()

# Step 3: paste the scratch sentence right in the middle of the rewritten
# run (right before the trailing ")."), which splits the original run into
# a prefix run / pasted run / suffix run with exactly the structure we want.
$insertAt = $runStart + "lignment Metamodels and Big Data. ".Length
$destination = $d.Range($insertAt, $insertAt)
$destination.Paste()

# Step 4: remove the scratch paragraph we used purely as a source for Copy.
$d.Paragraphs.Last.Range.Delete()

Write-Output "done"
